# Tableau_MVC.xlsx - "Début du formulaire de Pizza"
#
# The "lister_ingredients.php" action is replaced by a new
# "lister_compositions.php" action, a brand-new "afficher_form_composition.php"
# action/row is appended, the special highlight fill on the empty filler rows
# of the "Contrôleurs" sheet is removed, and the view is scrolled/selected on
# the new row.

$wb = $excel.ActiveWorkbook

$wsControleurs = $wb.Worksheets.Item("Contrôleurs")

# ---------------------------------------------------------------------------
# 1) Row 17 used to describe "lister_ingredients.php" - it now documents the
#    new "lister_compositions.php" action.
# ---------------------------------------------------------------------------
$wsControleurs.Range("C17").Value2 = "lister_compositions.php"
$wsControleurs.Range("A17").Value2 = "Retourne la liste de la composition de la pizza"

# Row 18 (finaliser_commande.php) keeps its text as-is.

# ---------------------------------------------------------------------------
# 2) Row 19 was an empty placeholder row - fill it with the new
#    "afficher_form_composition.php" action.
# ---------------------------------------------------------------------------
$wsControleurs.Range("A19").Value2 = "Récupérer le formulaire pour ajouter une composition"
$wsControleurs.Range("B19").Value2 = "JSON"
$wsControleurs.Range("C19").Value2 = "afficher_form_composition.php"
$wsControleurs.Range("D19").Value2 = "GET - type - Facultatif - Type de l'ingrédient pour filtrer le formulaire"

# ---------------------------------------------------------------------------
# 3) Rows 24-27 (the still-empty filler rows in column A) lose the special
#    peach/orange highlight fill they used to carry - copy the plain
#    "left + vertical-center" formatting (already used lower in the same
#    column, e.g. A6) onto them instead, which re-uses the existing style
#    instead of inventing a new one. (A23/A28/A29 already use that same
#    plain style, so nothing to do there.)
# ---------------------------------------------------------------------------
$wsControleurs.Range("A6").Copy() | Out-Null
$wsControleurs.Range("A24:A27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Move the view: scroll so row 10 is at the top and select the newly
#    filled D19 cell.
# ---------------------------------------------------------------------------
$wsControleurs.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$wsControleurs.Range("D19").Select()
